# Commit: "cambio formato datos de entrada (nueva columna con tipo de commodity)"
# Adds a "type" (commodity) column to comp_quantity_inst1 and two new
# aggregate rows (no_req_total / no_opt_total) on the parameters sheet.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("parameters")
$wsComp   = $wb.Worksheets.Item("comp_quantity_inst1")

# --- comp_quantity_inst1: new "type" column (commodity letter per row) ---
$wsComp.Range("E2").Value = "A"
$wsComp.Range("E3").Value = "B"
$wsComp.Range("E4").Value = "C"
$wsComp.Range("E5").Value = "D"
$wsComp.Range("E6").Value = "E"
$wsComp.Range("E7").Value = "F"
$wsComp.Range("E8").Value = "G"
$wsComp.Range("E9").Value = "H"
$wsComp.Range("E10").Value = "I"
$wsComp.Range("E11").Value = "J"
$wsComp.Range("E12").Value = "K"
$wsComp.Range("E1").Value = "type"

# whole used range now gets centered horizontal alignment
$wsComp.Range("A1:E12").HorizontalAlignment = -4108

# --- parameters: fleet_size updated + two new totals rows ---
$wsParams.Range("B12").Value = 22

$wsParams.Range("A13").Value = "no_req_total"
$wsParams.Range("B13").Formula = "=SUM(comp_quantity_inst1!C2:C7)"
$wsParams.Range("A14").Value = "no_opt_total"
$wsParams.Range("B14").Formula = "=SUM(comp_quantity_inst1!C8:C12)"

# --- view/selection state: parameters becomes the active tab ---
# (select on comp_quantity_inst1 first - selecting a range on it would
#  otherwise re-activate it and steal the "active tab" from parameters)
$wsComp.Range("C15").Select()

$wsParams.Activate()
$wsParams.Range("A13:A14").Select()
